$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (old D:K shifts to E:L) for the
# new FY2018 (period-ending 2018-12-31, serial 43465) financial data.
$ws.Columns("D").Insert()

# Copy number formats/styles from the (now-shifted) old column D, which
# landed in column E, into the newly inserted blank column D, per
# contiguous data block (skips label-only rows 5,6,37,79 which have no
# D:K cells).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Income Statement - Period Ending header
$ws.Range("D7").Value = 43465

# Income Statement data rows
$ws.Range("D8").Value = 1300
$ws.Range("D9").Value = 800
$ws.Range("D10").Value = 500
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 200
$ws.Range("D15").Value = 1200
$ws.Range("D17").Value = "NA"
$ws.Range("D18").Value = "NA"
$ws.Range("D20").Value = "NA"
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 500
$ws.Range("D23").Value = -5800
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = "NA"
$ws.Range("D27").Value = "NA"
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = "NA"
$ws.Range("D33").Value = "NA"
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = "NA"

# Balance Sheet - Period Ending header
$ws.Range("D38").Value = 43465

# Balance Sheet data rows
$ws.Range("D41").Value = 800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 500
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 200
$ws.Range("D46").Value = 1500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 36600
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 38100
$ws.Range("D57").Value = 700
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 1500
$ws.Range("D60").Value = 2200
$ws.Range("D61").Value = 17900
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 20100
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -89300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 18000
$ws.Range("D77").Value = 0

# Cash Flow Statement - Period Ending header
$ws.Range("D80").Value = 43465

# Cash Flow Statement data rows
$ws.Range("D81").Value = "NA"
$ws.Range("D83").Value = 1200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -1200
$ws.Range("D91").Value = -12100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -12100
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 13100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -200
